{"js": "// Change: \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c/\u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\" ->\n//   \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n// followed by a new list paragraph:\n//   \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst target = \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c/\u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\";\nlet found = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === target) {\n    found = paras.items[i];\n    break;\n  }\n}\n\nif (!found) {\n  throw new Error(\"Could not find paragraph: \" + target);\n}\n\n// Replace the whole paragraph's text (keeps its list/paragraph formatting).\nfound.insertText(\"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\", \"Replace\");\n\n// Insert a new list paragraph right after it with the split-off question.\nconst newPara = found.insertParagraph(\"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\", \"After\");\n\nawait context.sync();\n", "ps1": "# Change: \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c/\u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\" ->\n#   \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n# followed by a new list paragraph:\n#   \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c/\u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find target paragraph text\"\n}\n\n$rng = $find.Parent\n$rng.Text = \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u0438\u0437\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n\n# The found/updated range sits in the original list paragraph; grab it,\n# insert a fresh paragraph mark after it (inherits the same list formatting)\n# and fill the new paragraph with the split-off question.\n$para = $rng.Paragraphs(1)\n$para.Range.InsertParagraphAfter()\n$newPara = $para.Next()\n$newPara.Range.Text = \"\u041c\u043e\u0436\u0435\u0442 \u043b\u0438 \u043f\u043e\u0432\u0430\u0440 \u043e\u0442\u043c\u0435\u043d\u0438\u0442\u044c \u0437\u0430\u043a\u0430\u0437?\"\n"}
